$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)

# Row 2
$ws.Range('D2').Value = '30.661.56'
$ws.Range('E2').Value = '  -0.17%  '

# Row 3
$ws.Range('D3').Value = '1.915.83'
$ws.Range('E3').Value = '  +1.20%  '

# Row 4
$ws.Range('E4').Value = '  +0.13%  '

# Row 5
$ws.Range('D5').Value = '''239.56'
$ws.Range('E5').Value = '  -2.27%  '

# Row 6
$ws.Range('E6').Value = '  +0.10%  '

# Row 7
$ws.Range('D7').Value = '''0.4930'
$ws.Range('E7').Value = '  +0.15%  '

# Row 8
$ws.Range('D8').Value = '''0.2979'
$ws.Range('E8').Value = '  +0.70%  '

# Row 9
$ws.Range('D9').Value = '''0.06764'
$ws.Range('E9').Value = '  -0.45%  '

# Row 10
$ws.Range('D10').Value = '1.908.70'
$ws.Range('E10').Value = '  +1.16%  '

# Row 11
$ws.Range('D11').Value = '''17.17'
$ws.Range('E11').Value = '  -0.23%  '

# Row 12
$ws.Range('D12').Value = '''0.07360'
$ws.Range('E12').Value = '  +1.61%  '

# Row 13
$ws.Range('D13').Value = '''5.173'
$ws.Range('E13').Value = '  +2.62%  '

# Row 14
$ws.Range('D14').Value = '''88.83'
$ws.Range('E14').Value = '  -2.20%  '

# Row 15
$ws.Range('D15').Value = '''0.6708'
$ws.Range('E15').Value = '  -1.19%  '

# Row 16
$ws.Range('D16').Value = '30.646.31'
$ws.Range('E16').Value = '  -0.13%  '

# Row 17
$ws.Range('D17').Value = '''0.000007941'
$ws.Range('E17').Value = '  -0.55%  '

# Row 18
$ws.Range('E18').Value = '  +2.61%  '

# Row 19
$ws.Range('E19').Value = '  +0.03%  '

# Row 20
$ws.Range('B20').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C20').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D20').Value = '2.151.03'
$ws.Range('E20').Value = '  +0.92%  '

# Row 21
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '''5.343'
$ws.Range('E21').Value = '  +10.77%  '

# Row 22
$ws.Range('E22').Value = '  +0.19%  '

# Row 23
$ws.Range('D23').Value = '''203.47'
$ws.Range('E23').Value = '  +7.66%  '

# Row 24
$ws.Range('D24').Value = '''6.326'
$ws.Range('E24').Value = '  +2.90%  '

# Row 25
$ws.Range('D25').Value = '''9.642'
$ws.Range('E25').Value = '  +2.75%  '

# Row 26
$ws.Range('D26').Value = '''164.90'
$ws.Range('E26').Value = '  +5.73%  '

# Row 27
$ws.Range('D27').Value = '''18.84'
$ws.Range('E27').Value = '  -1.39%  '

# Row 28
$ws.Range('D28').Value = '''1.960'
$ws.Range('E28').Value = '  +3.07%  '

# Row 29
$ws.Range('D29').Value = '''1.480'
$ws.Range('E29').Value = '  +5.75%  '

# Row 30
$ws.Range('D30').Value = '''4.373'
$ws.Range('E30').Value = '  +0.76%  '

# Row 31
$ws.Range('D31').Value = '''0.09183'
$ws.Range('E31').Value = '  +1.17%  '

# Row 32
$ws.Range('D32').Value = '''4.064'
$ws.Range('E32').Value = '  +1.27%  '

# Row 33
$ws.Range('D33').Value = '''0.05274'
$ws.Range('E33').Value = '  +1.22%  '

# Row 34
$ws.Range('D34').Value = '''0.7425'
$ws.Range('E34').Value = '  -1.01%  '

# Row 35
$ws.Range('D35').Value = '''1.118'
$ws.Range('E35').Value = '  +0.85%  '

# Row 36
$ws.Range('D36').Value = '''2.727'
$ws.Range('E36').Value = '  -1.75%  '

# Row 37
$ws.Range('D37').Value = '''0.01848'
$ws.Range('E37').Value = '  +0.38%  '

# Row 38
$ws.Range('D38').Value = '''2.716'
$ws.Range('E38').Value = '  +1.05%  '

# Row 39
$ws.Range('D39').Value = '''0.9253'
$ws.Range('E39').Value = '  -1.33%  '

# Row 40
$ws.Range('D40').Value = '''2.078'
$ws.Range('E40').Value = '  -3.22%  '

# Row 41
$ws.Range('D41').Value = '''0.4464'
$ws.Range('E41').Value = '  +0.90%  '

# Row 42
$ws.Range('D42').Value = '''73.25'
$ws.Range('E42').Value = '  +27.00%  '

# Row 43
$ws.Range('D43').Value = '''5.990'
$ws.Range('E43').Value = '  +3.70%  '

# Row 44
$ws.Range('D44').Value = '''106.73'
$ws.Range('E44').Value = '  +1.27%  '

# Row 45
$ws.Range('E45').Value = '  +0.14%  '

# Row 46
$ws.Range('E46').Value = '  +3.80%  '

# Row 47
$ws.Range('D47').Value = '''7.644'
$ws.Range('E47').Value = '  +0.55%  '

# Row 48
$ws.Range('B48').Value = 'Elrond'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D48').Value = '''35.32'
$ws.Range('E48').Value = '  +5.02%  '

# Row 49
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '''9.022'
$ws.Range('E49').Value = '  +3.29%  '

# Row 50
$ws.Range('D50').Value = '''0.05883'
$ws.Range('E50').Value = '  +0.27%  '

# Row 51
$ws.Range('D51').Value = '''0.4036'
$ws.Range('E51').Value = '  +2.59%  '
